# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/string updates (Coin names, links, volume %, and multi-dot prices)
$ws.Range("D2").Value = '42.453.54'
$ws.Range("E2").Value = '  -2.46%  '
$ws.Range("D3").Value = '2.225.40'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E5").Value = '  -7.87%  '
$ws.Range("E6").Value = '  +10.84%  '
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("E9").Value = '  -2.39%  '
$ws.Range("E10").Value = '  -7.74%  '
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("E13").Value = '  -7.41%  '
$ws.Range("E14").Value = '  +14.94%  '
$ws.Range("E15").Value = '  -2.31%  '
$ws.Range("E16").Value = '  -4.39%  '
$ws.Range("D17").Value = '2.557.19'
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '2.229.04'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").Value = '42.490.45'
$ws.Range("E20").Value = '  +3.30%  '
$ws.Range("E21").Value = '  -3.58%  '
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("E23").Value = '  +19.60%  '
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("E25").Value = '  -2.16%  '
$ws.Range("E26").Value = '  -3.90%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E28").Value = '  -4.77%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E30").Value = '  -9.62%  '
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E32").Value = '  -5.67%  '
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("E36").Value = '  +11.31%  '
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("E38").Value = '  -2.29%  '
$ws.Range("E39").Value = '  -3.74%  '
$ws.Range("E40").Value = '  -5.42%  '
$ws.Range("E41").Value = '  +1.82%  '
$ws.Range("E42").Value = '  -5.57%  '
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  -9.66%  '
$ws.Range("E46").Value = '  -6.07%  '
$ws.Range("E47").Value = '  -6.37%  '
$ws.Range("E48").Value = '  +2.55%  '
$ws.Range("E49").Value = '  +10.54%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E51").Value = '  -0.53%  '

# Price values that would otherwise be auto-detected as numbers by Excel;
# force them to remain text to match the original inline-string formatting.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '294.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0915'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000105'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0885'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.13'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.127'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0372'
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.232'
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.54'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.71'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.47'
$ws.Range("D51").Style = "Normal"
